# Reorder the "Recorded By" contributor lists in column G for several rows.
# This only changes the ordering of the comma-separated names/emails within
# each cell; the set of values per row is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    3  = "System, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    4  = "servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    5  = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    6  = "mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"
    7  = "Fatmaelhady@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
    12 = "dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
    13 = "yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"
    20 = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    25 = "menna-allah.gamil@med.asu.edu.eg, Noran.Mahmoud@med.asu.edu.eg"
    30 = "aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
